# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" worksheets to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F on the "展览" sheet.
$updatesExhibit = @{
    7  = 577
    8  = 55
    9  = 8377
    12 = 1123
    13 = 888
    20 = 929
}

# Map of row number -> new value for column F on the "全部类型" sheet.
$updatesAll = @{
    9  = 577
    10 = 55
    11 = 8377
    14 = 1123
    15 = 888
    22 = 929
}

$wsExhibit = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibit.Keys) {
    $wsExhibit.Range("F$row").Value = $updatesExhibit[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $wsAll.Range("F$row").Value = $updatesAll[$row]
}
